$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.100.41"
$ws.Cells.Item(2, 5).Value = "  -0.50%  "
$ws.Cells.Item(3, 4).Value = "3.289.65"
$ws.Cells.Item(3, 5).Value = "  -0.64%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).Value = "'585.86"
$ws.Cells.Item(5, 5).Value = "  +2.59%  "
$ws.Cells.Item(6, 4).Value = "'180.21"
$ws.Cells.Item(6, 5).Value = "  -0.31%  "
$ws.Cells.Item(7, 4).Value = "'0.655"
$ws.Cells.Item(7, 5).Value = "  +9.45%  "
$ws.Cells.Item(8, 5).Value = "  +0.08%  "
$ws.Cells.Item(9, 4).Value = "'0.125"
$ws.Cells.Item(9, 5).Value = "  -2.62%  "
$ws.Cells.Item(10, 5).Value = "  +1.92%  "
$ws.Cells.Item(11, 4).Value = "'0.405"
$ws.Cells.Item(11, 5).Value = "  +0.65%  "
$ws.Cells.Item(12, 4).Value = "3.862.84"
$ws.Cells.Item(12, 5).Value = "  -0.48%  "
$ws.Cells.Item(13, 5).Value = "  -4.73%  "
$ws.Cells.Item(14, 4).Value = "66.127.90"
$ws.Cells.Item(14, 5).Value = "  -0.55%  "
$ws.Cells.Item(15, 4).Value = "'26.44"
$ws.Cells.Item(15, 5).Value = "  -2.12%  "
$ws.Cells.Item(16, 5).Value = "  -1.90%  "
$ws.Cells.Item(17, 4).Value = "3.295.78"
$ws.Cells.Item(17, 5).Value = "  +0.16%  "
$ws.Cells.Item(18, 4).Value = "'436.09"
$ws.Cells.Item(18, 5).Value = "  +0.06%  "
$ws.Cells.Item(19, 4).Value = "'13.25"
$ws.Cells.Item(19, 5).Value = "  -2.73%  "
$ws.Cells.Item(20, 5).Value = "  -3.07%  "
$ws.Cells.Item(21, 5).Value = "  -2.63%  "
$ws.Cells.Item(22, 4).Value = "'72.47"
$ws.Cells.Item(22, 5).Value = "  -1.24%  "
$ws.Cells.Item(23, 4).Value = "'0.999"
$ws.Cells.Item(23, 5).Value = "  -0.24%  "
$ws.Cells.Item(24, 5).Value = "  +0.22%  "
$ws.Cells.Item(25, 4).Value = "3.428.80"
$ws.Cells.Item(25, 5).Value = "  -0.56%  "
$ws.Cells.Item(26, 4).Value = "'0.510"
$ws.Cells.Item(26, 5).Value = "  -0.43%  "
$ws.Cells.Item(27, 5).Value = "  +4.12%  "
$ws.Cells.Item(28, 4).Value = "'0.0000113"
$ws.Cells.Item(28, 5).Value = "  -3.71%  "
$ws.Cells.Item(29, 5).Value = "  -1.49%  "
$ws.Cells.Item(30, 5).Value = "  +0.25%  "
$ws.Cells.Item(31, 5).Value = "  +0.73%  "
$ws.Cells.Item(32, 4).Value = "'22.33"
$ws.Cells.Item(32, 5).Value = "  -1.92%  "
$ws.Cells.Item(33, 5).Value = "  -0.02%  "
$ws.Cells.Item(34, 4).Value = "'5.20"
$ws.Cells.Item(34, 5).Value = "  -1.73%  "
$ws.Cells.Item(35, 4).Value = "'6.61"
$ws.Cells.Item(35, 5).Value = "  -2.17%  "
$ws.Cells.Item(36, 5).Value = "  -2.02%  "
$ws.Cells.Item(37, 4).Value = "'158.17"
$ws.Cells.Item(37, 5).Value = "  -0.82%  "
$ws.Cells.Item(38, 5).Value = "  -5.04%  "
$ws.Cells.Item(39, 4).Value = "'26.43"
$ws.Cells.Item(39, 5).Value = "  -3.02%  "
$ws.Cells.Item(40, 4).Value = "'1.77"
$ws.Cells.Item(40, 5).Value = "  -3.40%  "
$ws.Cells.Item(41, 4).Value = "2.795.83"
$ws.Cells.Item(41, 5).Value = "  +0.35%  "
$ws.Cells.Item(42, 4).Value = "'0.773"
$ws.Cells.Item(42, 5).Value = "  -1.41%  "
$ws.Cells.Item(43, 4).Value = "'4.34"
$ws.Cells.Item(43, 5).Value = "  -2.22%  "
$ws.Cells.Item(44, 4).Value = "'40.21"
$ws.Cells.Item(44, 5).Value = "  +0.19%  "
$ws.Cells.Item(45, 4).Value = "'6.07"
$ws.Cells.Item(45, 5).Value = "  -1.35%  "
$ws.Cells.Item(46, 5).Value = "  -1.76%  "
$ws.Cells.Item(47, 4).Value = "'2.31"
$ws.Cells.Item(47, 5).Value = "  -1.04%  "
$ws.Cells.Item(48, 4).Value = "'320.24"
$ws.Cells.Item(48, 5).Value = "  +0.25%  "
$ws.Cells.Item(49, 4).Value = "'23.21"
$ws.Cells.Item(49, 5).Value = "  -3.65%  "
$ws.Cells.Item(50, 4).Value = "'0.0267"
$ws.Cells.Item(50, 5).Value = "  -1.53%  "
$ws.Cells.Item(51, 5).Value = "  +7.18%  "
